$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.909.60"
$ws.Range("E2").Value = "  -2.61%  "
$ws.Range("D3").Value = "3.183.15"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.63%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.182.71"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.478"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.50%  "
$ws.Range("E13").Value = "  -3.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.58%  "
$ws.Range("D15").Value = "3.702.18"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "64.909.59"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "3.185.72"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "485.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.44%  "
$ws.Range("E29").Value = "  -4.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.76%  "
$ws.Range("E32").Value = "  -7.65%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -5.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.77%  "
$ws.Range("D39").Value = "0.0₃0741"
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "463.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.128"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0407"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "2.928.79"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.278"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.31%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.117"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.36%  "
